$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "Day_of_experiment"

# Day_of_experiment values for rows 2-21 (five rows of 1, five rows of 2, repeating
# for each Cell_type block)
$dayValues = @(1,1,1,1,1,2,2,2,2,2,1,1,1,1,1,2,2,2,2,2)

for ($i = 0; $i -lt $dayValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $dayValues[$i]
}

# Best-fit the new column, matching the width Excel picked for the header text
$ws.Columns.Item(5).ColumnWidth = 16.498697916666668

# Move the active selection to the new bottom-right corner of the data, as Excel
# would leave it after typing the last value
$ws.Range("E21").Select()
